# Updates cryptos list values (price/volume columns) per the scraped commit diff.
# Rows 13/14 and 46/47 also have their Coin/Link swapped (ranking re-sort).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '76.378.06'
$ws.Range("E2").Value = '  +0.54%  '
# Row 3
$ws.Range("D3").Value = '2.954.85'
$ws.Range("E3").Value = '  +2.25%  '
# Row 4
$ws.Range("E4").Value = '  +0.06%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '199.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.85%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '596.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.16%  '
# Row 7
$ws.Range("E7").Value = '  +0.03%  '
# Row 8
$ws.Range("E8").Value = '  -0.29%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.201'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.43%  '
# Row 10
$ws.Range("D10").Value = '2.956.37'
$ws.Range("E10").Value = '  +2.30%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.444'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.37%  '
# Row 12
$ws.Range("E12").Value = '  +0.47%  '
# Row 13
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.42%  '
# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.501.96'
$ws.Range("E14").Value = '  +4.48%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.31%  '
# Row 16
$ws.Range("D16").Value = '76.363.34'
$ws.Range("E16").Value = '  +0.67%  '
# Row 17
$ws.Range("E17").Value = '  -0.43%  '
# Row 18
$ws.Range("D18").Value = '2.950.10'
$ws.Range("E18").Value = '  +2.63%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.65%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.78'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.19%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.14%  '
# Row 23
$ws.Range("E23").Value = '  +4.18%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.74%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '
# Row 26
$ws.Range("D26").Value = '3.106.17'
$ws.Range("E26").Value = '  +2.72%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.58%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.72'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.10%  '
# Row 29
$ws.Range("E29").Value = '  +0.22%  '
# Row 30
$ws.Range("E30").Value = '  -0.02%  '
# Row 31
$ws.Range("E31").Value = '  +10.27%  '
# Row 32
$ws.Range("E32").Value = '  -1.69%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '497.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.43%  '
# Row 34
$ws.Range("E34").Value = '  +0.90%  '
# Row 35
$ws.Range("E35").Value = '  +0.05%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.57%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.57%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.392'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +13.46%  '
# Row 39
$ws.Range("E39").Value = '  +18.35%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.32%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.111'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.32%  '
# Row 42
$ws.Range("E42").Value = '  +0.02%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '180.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.93%  '
# Row 44
$ws.Range("E44").Value = '  -1.32%  '
# Row 45
$ws.Range("E45").Value = '  -1.70%  '
# Row 46
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.30%  '
# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.69%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.591'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.03%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.71%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.33'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.95%  '
# Row 51
$ws.Range("E51").Value = '  +0.41%  '
